# Insert a new data row at row 23 (pushing existing rows 23-106 down to 24-107),
# matching the weekly update described in the commit message
# ("Fruta / hortaliza, semanal" -> a new week's record was added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23; this shifts rows 23..106 down to 24..107
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record's data.
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44971
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112022
$ws.Range("G23").Value = "Arveja Verde"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("M23").Value = 22000
$ws.Range("N23").Value = "`$/saco 25 kilos"
$ws.Range("O23").Value = "Provincia de Diguillín"
$ws.Range("P23").Value = 880
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
